$d = $word.ActiveDocument

# --- Change 1: simplify the context-diagram sentence, removing
#     "та менеджера катологу" (the misspelled "catalog manager" actor). ---
$d.Content.Find.Execute(
    "адміністратора та менеджера катологу. ",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "адміністратора. ",
    2)

# --- Change 2: move the "_GoBack" bookmark from the end of the API
#     endpoint listing to the end of the "...із відповідними сервісами."
#     paragraph. ---

# Remove the old bookmark (currently sitting after the closing "}" of
# the last API-endpoint line).
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# Locate the end of the target sentence.
$target = $d.Content.Duplicate
$target.Find.Execute(
    "із відповідними сервісами.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "",
    0)
$target.Collapse(0)

# Inserting a bookmark directly at a collapsed position sitting right at
# a paragraph-mark boundary lands in the wrong spot, so insert a
# throwaway character, bookmark it, then remove the character again --
# this leaves the (now collapsed) bookmark correctly placed right after
# the sentence and before the paragraph mark.
$target.InsertAfter("X")
$marker = $d.Range($target.Start, $target.Start + 1)
$d.Bookmarks.Add("_GoBack", $marker)
$marker.Text = ""
